$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph by scanning the paragraph
# collection (robust against any Find()-returned Range quirks).
$count = $d.Paragraphs.Count
$jupiterIdx = 0
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("Ver no Jupiter")) {
        $jupiterIdx = $i
    }
}

if ($jupiterIdx -eq 0) {
    throw "Could not find the 'Ver no Jupiter' paragraph"
}

# Remove the three paragraphs that followed the "LOQ4237: ..." requirement
# paragraph: the blank spacer paragraph right before it, the "Ver no
# Jupiter ..." paragraph itself, and the "(c) 2020 ... Creative Commons
# Attribution" paragraph right after it. The blank paragraph that sits
# right before the page-break paragraph is left untouched.
$startPara = $d.Paragraphs.Item($jupiterIdx - 1)
$endPara = $d.Paragraphs.Item($jupiterIdx + 1)

$r = $d.Range($startPara.Range.Start, $endPara.Range.End)
$r.Delete()
